$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-crawl timestamp: every data row's "timestamp" column (O) moves
#    from 2022-07-31 07:00:59 to 2022-07-31 20:58:28 (rows 2..130).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 130; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-07-31 20:58:28"
}

# ---------------------------------------------------------------------
# 2) Two products swapped list position: the "Avela Soeckchen" item
#    (row 27) and the "Avela Strumpfhose" item (row 28) traded places.
#    Capture the original column A..N content/flags for both rows,
#    then write them back swapped. Columns holding numeric-looking
#    text (id, price, priceContextPrice) are re-applied as Text so
#    Excel doesn't silently re-type them as numbers.
# ---------------------------------------------------------------------
$cols = 1..14
# Columns whose values are digit/decimal strings that must stay text.
$textCols = @(1, 8, 11)

$row27 = @{}
$row28 = @{}
foreach ($c in $cols) {
    $row27[$c] = $ws.Cells.Item(27, $c).Value2
    $row28[$c] = $ws.Cells.Item(28, $c).Value2
}

foreach ($c in $cols) {
    $target27 = $row28[$c]
    $target28 = $row27[$c]

    $cell27 = $ws.Cells.Item(27, $c)
    $cell28 = $ws.Cells.Item(28, $c)

    if ($textCols -contains $c) {
        if ($target27 -ne "") {
            $cell27.NumberFormat = "@"
        }
        if ($target28 -ne "") {
            $cell28.NumberFormat = "@"
        }
    }

    $cell27.Value = $target27
    $cell28.Value = $target28

    if ($textCols -contains $c) {
        if ($target27 -ne "") {
            $cell27.Style = "Normal"
        }
        if ($target28 -ne "") {
            $cell28.Style = "Normal"
        }
    }
}
